# Actualización automática hashcode lun jul  8 02:08:00 CEST 2019
# Updates the "hashcode" column (column B) values for a set of rows
# identified by their matching key in column A, to new hash values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B44"  = "775da89266fde57dfe7ca7c89abf5d91"
    "B74"  = "8a74666dc4ebb183229cedc771aa374f"
    "B89"  = "e5a9c26e094a5557ae9c4aa83e416d55"
    "B99"  = "0c473cacc596f7b80f753639d0d0ca9c"
    "B110" = "8c9098805d070995ea6995c660cc73a1"
    "B154" = "e9828e955ed4896624069e2230da5da2"
    "B160" = "f3de5288eeaf606f566c40f38f1f948a"
    "B161" = "1e5c3f3bf56fea72588394470e1cc359"
    "B168" = "bc95cae257a5ff8399d8aa38ac0096e0"
    "B278" = "9283cf6e227051ed64790cd8214746ac"
    "B330" = "02d08555a89aca4227289c60c19d9b82"
    "B345" = "3d3502f758d76be92c0f4e2ea3201dd1"
    "B534" = "76da3783aa2a61aa6867b6ba825b3179"
    "B547" = "61c4f18193adac7d146bc75c0f680430"
    "B553" = "58d85ba2051dd71507a5e4255d2e5b94"
    "B768" = "856d009b685edcaa25e7aebd1e4cb92c"
    "B811" = "dbd952bba9bedbb15ced3d14a76bc9b0"
    "B815" = "bd5b9380588c9dc7c9ba8123dc3cab76"
    "B816" = "831b12f239db1883cfb6a62cd480eabe"
    "B825" = "e0b748b7abab51601ff88878e1646e1d"
    "B827" = "e72e4ad52475855fd285dd2b5bbecbd4"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
